# Generate Report for Handoff
# Updates the localization-status workbook to reflect that "b.md" has now
# been handed off (zh-cn fully; de-de as well), while the existing
# zh-cn status for "a.md" is normalized to "Ready for handoff" too, and a
# version-mismatch error detail is recorded for the new b.md handoff.

$wb = $excel.ActiveWorkbook

$newHandoffDate_zhcn = "2016-08-27 10:38:07"
$newHandoffDate_dede = "2016-08-27 10:38:11"
$newHandoffFile_zhcn = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$newHandoffFile_dede = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$readyStatus = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4323094c55b0ef4dc6d7ebdbfb31a2f5c1368cb4/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/408147debbc25497a45f46e86ef886d21bdb2262/e2e/b.md."

# ---- Overview sheet: row 3 is b.md ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $readyStatus
$wsOverview.Range("F3").Value = $readyStatus
$wsOverview.Range("G3").Value = $newHandoffDate_dede

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# row 2 = a.md: status normalized
$wsZhCn.Range("C2").Value = $readyStatus
# row 3 = b.md: new handoff round
$wsZhCn.Range("C3").Value = $readyStatus
$wsZhCn.Range("G3").Value = $newHandoffFile_zhcn
$wsZhCn.Range("H3").Value = $newHandoffDate_zhcn
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
# row 3 = b.md: new handoff round
$wsDeDe.Range("C3").Value = $readyStatus
$wsDeDe.Range("G3").Value = $newHandoffFile_dede
$wsDeDe.Range("H3").Value = $newHandoffDate_dede
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
